$wb = $excel.ActiveWorkbook
$after = $wb.Worksheets.Item("recording steps and topics")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$newSheet.Name = "Test"
$newSheet.Columns.Item(1).ColumnWidth = 15.57
Write-Host "Done"
